# Homework and lab updates.
#
# The worksheet tab labelled "Meetups" is the one that actually contains the
# weekly schedule data (Date / Topic / .../ Homework / Lab columns) - the
# workbook's internal file naming and tab naming are swapped relative to
# what you might expect, but "Meetups" is the sheet with the Homework/Lab
# link columns that this change targets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Row 10 (Multiple Regression): replace the generic homework9.pdf link with
# the more specific "multiple regression" homework file. Set this first so
# that the new shared string is appended in the same order as the source
# workbook (homework9_mr.pdf, then homework9_lr.pdf, then the new lab zip).
$ws.Range("G10").Value = "/static/homework/homework9_mr.pdf"

# Row 9 (Maximum Likelihood Estimation and Logistic Regression): these cells
# were previously blank placeholders (G9 only had formatting, H9 had no cell
# at all). Clear the leftover placeholder style on G9 before writing so it
# matches a plain text cell like its neighbors, then fill in the new
# homework and lab links.
$ws.Range("G9").Style = "Normal"
$ws.Range("G9").Value = "/static/homework/homework9_lr.pdf"
$ws.Range("H9").Value = "/static/labs/09_logistic_regression.zip"

# Update the saved cell selection/cursor position to match the edited area.
$ws.Range("H10").Select() | Out-Null
